$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44257
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 6000
$ws.Cells.Item(2, 15).Value = 6500
$ws.Cells.Item(2, 16).Value = 6250
$ws.Cells.Item(2, 18).Value = 'Perú'
$ws.Cells.Item(2, 19).Value = 1562

$ws.Cells.Item(3, 4).Value = 44363
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 6500
$ws.Cells.Item(3, 15).Value = 7000
$ws.Cells.Item(3, 16).Value = 6750
$ws.Cells.Item(3, 18).Value = 'Perú'
$ws.Cells.Item(3, 19).Value = 1688

$ws.Cells.Item(4, 4).Value = 44581
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 150
$ws.Cells.Item(4, 14).Value = 7000
$ws.Cells.Item(4, 15).Value = 7000
$ws.Cells.Item(4, 16).Value = 7000
$ws.Cells.Item(4, 18).Value = 'Perú'
$ws.Cells.Item(4, 19).Value = 1750

$ws.Cells.Item(5, 4).Value = 44462
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 7500
$ws.Cells.Item(5, 15).Value = 8000
$ws.Cells.Item(5, 16).Value = 7750
$ws.Cells.Item(5, 18).Value = 'Brasil'
$ws.Cells.Item(5, 19).Value = 1938

$ws.Cells.Item(6, 4).Value = 44442
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 200
$ws.Cells.Item(6, 14).Value = 8000
$ws.Cells.Item(6, 15).Value = 8500
$ws.Cells.Item(6, 16).Value = 8250
$ws.Cells.Item(6, 18).Value = 'Brasil'
$ws.Cells.Item(6, 19).Value = 2062

$ws.Cells.Item(7, 4).Value = 44398
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 7500
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 7750
$ws.Cells.Item(7, 18).Value = 'Brasil'
$ws.Cells.Item(7, 19).Value = 1938

$ws.Cells.Item(8, 4).Value = 44588
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 6500
$ws.Cells.Item(8, 15).Value = 7000
$ws.Cells.Item(8, 16).Value = 6800
$ws.Cells.Item(8, 18).Value = 'Perú'
$ws.Cells.Item(8, 19).Value = 1700

$ws.Cells.Item(9, 4).Value = 44335
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 9500
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 9750
$ws.Cells.Item(9, 18).Value = 'Perú'
$ws.Cells.Item(9, 19).Value = 2438

$ws.Cells.Item(10, 4).Value = 44671
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 7500
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 7750
$ws.Cells.Item(10, 18).Value = 'Ecuador'
$ws.Cells.Item(10, 19).Value = 1938

$ws.Cells.Item(11, 4).Value = 44477
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 200
$ws.Cells.Item(11, 14).Value = 8000
$ws.Cells.Item(11, 15).Value = 8500
$ws.Cells.Item(11, 16).Value = 8250
$ws.Cells.Item(11, 18).Value = 'Perú'
$ws.Cells.Item(11, 19).Value = 2062

$ws.Cells.Item(12, 4).Value = 44316
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 9000
$ws.Cells.Item(12, 15).Value = 10000
$ws.Cells.Item(12, 16).Value = 9500
$ws.Cells.Item(12, 18).Value = 'Ecuador'
$ws.Cells.Item(12, 19).Value = 2375

$ws.Cells.Item(13, 4).Value = 44526
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 170
$ws.Cells.Item(13, 14).Value = 6000
$ws.Cells.Item(13, 15).Value = 6500
$ws.Cells.Item(13, 16).Value = 6235
$ws.Cells.Item(13, 18).Value = 'Perú'
$ws.Cells.Item(13, 19).Value = 1559

$ws.Cells.Item(14, 4).Value = 44574
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 6500
$ws.Cells.Item(14, 15).Value = 7000
$ws.Cells.Item(14, 16).Value = 6750
$ws.Cells.Item(14, 18).Value = 'Perú'
$ws.Cells.Item(14, 19).Value = 1688

$ws.Cells.Item(15, 4).Value = 44446
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 7000
$ws.Cells.Item(15, 15).Value = 7500
$ws.Cells.Item(15, 16).Value = 7250
$ws.Cells.Item(15, 18).Value = 'Brasil'
$ws.Cells.Item(15, 19).Value = 1812

$ws.Cells.Item(16, 4).Value = 44405
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 8000
$ws.Cells.Item(16, 15).Value = 8500
$ws.Cells.Item(16, 16).Value = 8250
$ws.Cells.Item(16, 18).Value = 'Brasil'
$ws.Cells.Item(16, 19).Value = 2062

$ws.Cells.Item(17, 4).Value = 44208
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 6000
$ws.Cells.Item(17, 15).Value = 6500
$ws.Cells.Item(17, 16).Value = 6250
$ws.Cells.Item(17, 18).Value = 'Perú'
$ws.Cells.Item(17, 19).Value = 1562

$ws.Cells.Item(18, 4).Value = 44474
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 8500
$ws.Cells.Item(18, 15).Value = 9000
$ws.Cells.Item(18, 16).Value = 8750
$ws.Cells.Item(18, 18).Value = 'Perú'
$ws.Cells.Item(18, 19).Value = 2188

$ws.Cells.Item(19, 4).Value = 44610
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 180
$ws.Cells.Item(19, 14).Value = 7500
$ws.Cells.Item(19, 15).Value = 8000
$ws.Cells.Item(19, 16).Value = 7722
$ws.Cells.Item(19, 18).Value = 'Perú'
$ws.Cells.Item(19, 19).Value = 1930

$ws.Cells.Item(20, 4).Value = 44211
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(20, 14).Value = 6000
$ws.Cells.Item(20, 15).Value = 6500
$ws.Cells.Item(20, 16).Value = 6250
$ws.Cells.Item(20, 18).Value = 'Perú'
$ws.Cells.Item(20, 19).Value = 1562

$ws.Cells.Item(21, 4).Value = 44490
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 220
$ws.Cells.Item(21, 14).Value = 6500
$ws.Cells.Item(21, 15).Value = 7000
$ws.Cells.Item(21, 16).Value = 6727
$ws.Cells.Item(21, 18).Value = 'Perú'
$ws.Cells.Item(21, 19).Value = 1682

$ws.Cells.Item(22, 4).Value = 44188
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 5500
$ws.Cells.Item(22, 15).Value = 6000
$ws.Cells.Item(22, 16).Value = 5750
$ws.Cells.Item(22, 18).Value = 'Perú'
$ws.Cells.Item(22, 19).Value = 1438

$ws.Cells.Item(23, 4).Value = 44630
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 200
$ws.Cells.Item(23, 14).Value = 7000
$ws.Cells.Item(23, 15).Value = 7500
$ws.Cells.Item(23, 16).Value = 7250
$ws.Cells.Item(23, 18).Value = 'Perú'
$ws.Cells.Item(23, 19).Value = 1812

$ws.Cells.Item(24, 4).Value = 44299
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 400
$ws.Cells.Item(24, 14).Value = 6000
$ws.Cells.Item(24, 15).Value = 6500
$ws.Cells.Item(24, 16).Value = 6250
$ws.Cells.Item(24, 18).Value = 'Perú'
$ws.Cells.Item(24, 19).Value = 1562

$ws.Cells.Item(25, 4).Value = 44575
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 200
$ws.Cells.Item(25, 14).Value = 6500
$ws.Cells.Item(25, 15).Value = 7000
$ws.Cells.Item(25, 16).Value = 6750
$ws.Cells.Item(25, 18).Value = 'Ecuador'
$ws.Cells.Item(25, 19).Value = 1688

$ws.Cells.Item(26, 4).Value = 44322
$ws.Cells.Item(26, 12).Value = 'Primera'
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 9000
$ws.Cells.Item(26, 15).Value = 10000
$ws.Cells.Item(26, 16).Value = 9500
$ws.Cells.Item(26, 18).Value = 'Perú'
$ws.Cells.Item(26, 19).Value = 2375

$ws.Cells.Item(27, 4).Value = 44715
$ws.Cells.Item(27, 12).Value = 'Primera'
$ws.Cells.Item(27, 13).Value = 300
$ws.Cells.Item(27, 14).Value = 9000
$ws.Cells.Item(27, 15).Value = 10000
$ws.Cells.Item(27, 16).Value = 9333
$ws.Cells.Item(27, 18).Value = 'Ecuador'
$ws.Cells.Item(27, 19).Value = 2333

$ws.Cells.Item(28, 4).Value = 44740
$ws.Cells.Item(28, 12).Value = 'Primera'
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 8000
$ws.Cells.Item(28, 15).Value = 8500
$ws.Cells.Item(28, 16).Value = 8250
$ws.Cells.Item(28, 18).Value = 'Brasil'
$ws.Cells.Item(28, 19).Value = 2062

$ws.Cells.Item(29, 4).Value = 44166
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 7000
$ws.Cells.Item(29, 15).Value = 7500
$ws.Cells.Item(29, 16).Value = 7250
$ws.Cells.Item(29, 18).Value = 'Perú'
$ws.Cells.Item(29, 19).Value = 1812

$ws.Cells.Item(30, 4).Value = 44559
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 300
$ws.Cells.Item(30, 14).Value = 6000
$ws.Cells.Item(30, 15).Value = 6500
$ws.Cells.Item(30, 16).Value = 6333
$ws.Cells.Item(30, 18).Value = 'Perú'
$ws.Cells.Item(30, 19).Value = 1583

$ws.Cells.Item(31, 4).Value = 44309
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 200
$ws.Cells.Item(31, 14).Value = 9500
$ws.Cells.Item(31, 15).Value = 10000
$ws.Cells.Item(31, 16).Value = 9750
$ws.Cells.Item(31, 18).Value = 'Perú'
$ws.Cells.Item(31, 19).Value = 2438

$ws.Cells.Item(32, 4).Value = 44742
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 200
$ws.Cells.Item(32, 14).Value = 8000
$ws.Cells.Item(32, 15).Value = 8500
$ws.Cells.Item(32, 16).Value = 8250
$ws.Cells.Item(32, 18).Value = 'Brasil'
$ws.Cells.Item(32, 19).Value = 2062

$ws.Cells.Item(33, 4).Value = 44294
$ws.Cells.Item(33, 12).Value = 'Primera'
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(33, 14).Value = 6500
$ws.Cells.Item(33, 15).Value = 7000
$ws.Cells.Item(33, 16).Value = 6750
$ws.Cells.Item(33, 18).Value = 'Perú'
$ws.Cells.Item(33, 19).Value = 1688

$ws.Cells.Item(34, 4).Value = 44608
$ws.Cells.Item(34, 12).Value = 'Primera'
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 7000
$ws.Cells.Item(34, 15).Value = 8000
$ws.Cells.Item(34, 16).Value = 7500
$ws.Cells.Item(34, 18).Value = 'Perú'
$ws.Cells.Item(34, 19).Value = 1875

$ws.Cells.Item(35, 4).Value = 44320
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 10000
$ws.Cells.Item(35, 15).Value = 11000
$ws.Cells.Item(35, 16).Value = 10500
$ws.Cells.Item(35, 18).Value = 'Perú'
$ws.Cells.Item(35, 19).Value = 2625

$ws.Cells.Item(36, 4).Value = 44246
$ws.Cells.Item(36, 12).Value = 'Primera'
$ws.Cells.Item(36, 13).Value = 400
$ws.Cells.Item(36, 14).Value = 6000
$ws.Cells.Item(36, 15).Value = 6500
$ws.Cells.Item(36, 16).Value = 6250
$ws.Cells.Item(36, 18).Value = 'Perú'
$ws.Cells.Item(36, 19).Value = 1562

$ws.Cells.Item(37, 4).Value = 44602
$ws.Cells.Item(37, 12).Value = 'Primera'
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 6500
$ws.Cells.Item(37, 15).Value = 7000
$ws.Cells.Item(37, 16).Value = 6750
$ws.Cells.Item(37, 18).Value = 'Perú'
$ws.Cells.Item(37, 19).Value = 1688

$ws.Cells.Item(38, 4).Value = 44467
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 300
$ws.Cells.Item(38, 14).Value = 7500
$ws.Cells.Item(38, 15).Value = 8000
$ws.Cells.Item(38, 16).Value = 7667
$ws.Cells.Item(38, 18).Value = 'Brasil'
$ws.Cells.Item(38, 19).Value = 1917

$ws.Cells.Item(39, 4).Value = 44516
$ws.Cells.Item(39, 12).Value = 'Primera'
$ws.Cells.Item(39, 13).Value = 400
$ws.Cells.Item(39, 14).Value = 7500
$ws.Cells.Item(39, 15).Value = 8000
$ws.Cells.Item(39, 16).Value = 7750
$ws.Cells.Item(39, 18).Value = 'Perú'
$ws.Cells.Item(39, 19).Value = 1938

$ws.Cells.Item(40, 4).Value = 44281
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 400
$ws.Cells.Item(40, 14).Value = 6000
$ws.Cells.Item(40, 15).Value = 6500
$ws.Cells.Item(40, 16).Value = 6250
$ws.Cells.Item(40, 18).Value = 'Perú'
$ws.Cells.Item(40, 19).Value = 1562

$ws.Cells.Item(41, 4).Value = 44714
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 200
$ws.Cells.Item(41, 14).Value = 10000
$ws.Cells.Item(41, 15).Value = 11000
$ws.Cells.Item(41, 16).Value = 10500
$ws.Cells.Item(41, 18).Value = 'Ecuador'
$ws.Cells.Item(41, 19).Value = 2625

$ws.Cells.Item(42, 4).Value = 44530
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 200
$ws.Cells.Item(42, 14).Value = 7500
$ws.Cells.Item(42, 15).Value = 8000
$ws.Cells.Item(42, 16).Value = 7750
$ws.Cells.Item(42, 18).Value = 'Perú'
$ws.Cells.Item(42, 19).Value = 1938

$ws.Cells.Item(43, 4).Value = 44272
$ws.Cells.Item(43, 12).Value = 'Primera'
$ws.Cells.Item(43, 13).Value = 200
$ws.Cells.Item(43, 14).Value = 6000
$ws.Cells.Item(43, 15).Value = 6500
$ws.Cells.Item(43, 16).Value = 6250
$ws.Cells.Item(43, 18).Value = 'Perú'
$ws.Cells.Item(43, 19).Value = 1562

$ws.Cells.Item(44, 4).Value = 44524
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 150
$ws.Cells.Item(44, 14).Value = 6500
$ws.Cells.Item(44, 15).Value = 7000
$ws.Cells.Item(44, 16).Value = 6667
$ws.Cells.Item(44, 18).Value = 'Perú'
$ws.Cells.Item(44, 19).Value = 1667

$ws.Cells.Item(45, 4).Value = 44699
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 200
$ws.Cells.Item(45, 14).Value = 7000
$ws.Cells.Item(45, 15).Value = 8000
$ws.Cells.Item(45, 16).Value = 7500
$ws.Cells.Item(45, 18).Value = 'Brasil'
$ws.Cells.Item(45, 19).Value = 1875

$ws.Cells.Item(46, 4).Value = 44350
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 200
$ws.Cells.Item(46, 14).Value = 9000
$ws.Cells.Item(46, 15).Value = 10000
$ws.Cells.Item(46, 16).Value = 9500
$ws.Cells.Item(46, 18).Value = 'Perú'
$ws.Cells.Item(46, 19).Value = 2375

$ws.Cells.Item(47, 4).Value = 44250
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 400
$ws.Cells.Item(47, 14).Value = 6000
$ws.Cells.Item(47, 15).Value = 6500
$ws.Cells.Item(47, 16).Value = 6250
$ws.Cells.Item(47, 18).Value = 'Perú'
$ws.Cells.Item(47, 19).Value = 1562

$ws.Cells.Item(48, 4).Value = 44245
$ws.Cells.Item(48, 12).Value = 'Primera'
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 6000
$ws.Cells.Item(48, 15).Value = 6500
$ws.Cells.Item(48, 16).Value = 6250
$ws.Cells.Item(48, 18).Value = 'Perú'
$ws.Cells.Item(48, 19).Value = 1562

$ws.Cells.Item(49, 4).Value = 44264
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 400
$ws.Cells.Item(49, 14).Value = 6000
$ws.Cells.Item(49, 15).Value = 6500
$ws.Cells.Item(49, 16).Value = 6250
$ws.Cells.Item(49, 18).Value = 'Perú'
$ws.Cells.Item(49, 19).Value = 1562

$ws.Cells.Item(50, 4).Value = 44159
$ws.Cells.Item(50, 12).Value = 'Primera'
$ws.Cells.Item(50, 13).Value = 400
$ws.Cells.Item(50, 14).Value = 6500
$ws.Cells.Item(50, 15).Value = 7000
$ws.Cells.Item(50, 16).Value = 6750
$ws.Cells.Item(50, 18).Value = 'Perú'
$ws.Cells.Item(50, 19).Value = 1688

$ws.Cells.Item(51, 4).Value = 44168
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 7000
$ws.Cells.Item(51, 15).Value = 7500
$ws.Cells.Item(51, 16).Value = 7250
$ws.Cells.Item(51, 18).Value = 'Perú'
$ws.Cells.Item(51, 19).Value = 1812

$ws.Cells.Item(52, 4).Value = 44232
$ws.Cells.Item(52, 12).Value = 'Primera'
$ws.Cells.Item(52, 13).Value = 400
$ws.Cells.Item(52, 14).Value = 6000
$ws.Cells.Item(52, 15).Value = 6500
$ws.Cells.Item(52, 16).Value = 6250
$ws.Cells.Item(52, 18).Value = 'Perú'
$ws.Cells.Item(52, 19).Value = 1562

$ws.Cells.Item(53, 4).Value = 44225
$ws.Cells.Item(53, 12).Value = 'Primera'
$ws.Cells.Item(53, 13).Value = 400
$ws.Cells.Item(53, 14).Value = 6000
$ws.Cells.Item(53, 15).Value = 6500
$ws.Cells.Item(53, 16).Value = 6250
$ws.Cells.Item(53, 18).Value = 'Ecuador'
$ws.Cells.Item(53, 19).Value = 1562

$ws.Cells.Item(54, 4).Value = 44587
$ws.Cells.Item(54, 12).Value = 'Primera'
$ws.Cells.Item(54, 13).Value = 200
$ws.Cells.Item(54, 14).Value = 7000
$ws.Cells.Item(54, 15).Value = 7500
$ws.Cells.Item(54, 16).Value = 7250
$ws.Cells.Item(54, 18).Value = 'Perú'
$ws.Cells.Item(54, 19).Value = 1812

$ws.Cells.Item(55, 4).Value = 44665
$ws.Cells.Item(55, 12).Value = 'Primera'
$ws.Cells.Item(55, 13).Value = 180
$ws.Cells.Item(55, 14).Value = 6500
$ws.Cells.Item(55, 15).Value = 7000
$ws.Cells.Item(55, 16).Value = 6778
$ws.Cells.Item(55, 18).Value = 'Perú'
$ws.Cells.Item(55, 19).Value = 1694

$ws.Cells.Item(56, 4).Value = 44217
$ws.Cells.Item(56, 12).Value = 'Primera'
$ws.Cells.Item(56, 13).Value = 200
$ws.Cells.Item(56, 14).Value = 6000
$ws.Cells.Item(56, 15).Value = 6500
$ws.Cells.Item(56, 16).Value = 6250
$ws.Cells.Item(56, 18).Value = 'Perú'
$ws.Cells.Item(56, 19).Value = 1562

$ws.Cells.Item(57, 4).Value = 44505
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 200
$ws.Cells.Item(57, 14).Value = 7000
$ws.Cells.Item(57, 15).Value = 7500
$ws.Cells.Item(57, 16).Value = 7250
$ws.Cells.Item(57, 18).Value = 'Perú'
$ws.Cells.Item(57, 19).Value = 1812

$ws.Cells.Item(58, 4).Value = 44505
$ws.Cells.Item(58, 12).Value = 'Segunda'
$ws.Cells.Item(58, 13).Value = 100
$ws.Cells.Item(58, 14).Value = 6500
$ws.Cells.Item(58, 15).Value = 6500
$ws.Cells.Item(58, 16).Value = 6500
$ws.Cells.Item(58, 18).Value = 'Perú'
$ws.Cells.Item(58, 19).Value = 1625

$ws.Cells.Item(59, 4).Value = 44329
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 200
$ws.Cells.Item(59, 14).Value = 10000
$ws.Cells.Item(59, 15).Value = 11000
$ws.Cells.Item(59, 16).Value = 10500
$ws.Cells.Item(59, 18).Value = 'Perú'
$ws.Cells.Item(59, 19).Value = 2625

$ws.Cells.Item(60, 4).Value = 44306
$ws.Cells.Item(60, 12).Value = 'Primera'
$ws.Cells.Item(60, 13).Value = 400
$ws.Cells.Item(60, 14).Value = 10000
$ws.Cells.Item(60, 15).Value = 11000
$ws.Cells.Item(60, 16).Value = 10500
$ws.Cells.Item(60, 18).Value = 'Perú'
$ws.Cells.Item(60, 19).Value = 2625

$ws.Cells.Item(61, 4).Value = 44435
$ws.Cells.Item(61, 12).Value = 'Primera'
$ws.Cells.Item(61, 13).Value = 200
$ws.Cells.Item(61, 14).Value = 8000
$ws.Cells.Item(61, 15).Value = 8500
$ws.Cells.Item(61, 16).Value = 8250
$ws.Cells.Item(61, 18).Value = 'Brasil'
$ws.Cells.Item(61, 19).Value = 2062

$ws.Cells.Item(62, 4).Value = 44519
$ws.Cells.Item(62, 12).Value = 'Primera'
$ws.Cells.Item(62, 13).Value = 200
$ws.Cells.Item(62, 14).Value = 7000
$ws.Cells.Item(62, 15).Value = 7500
$ws.Cells.Item(62, 16).Value = 7250
$ws.Cells.Item(62, 18).Value = 'Perú'
$ws.Cells.Item(62, 19).Value = 1812

$ws.Cells.Item(63, 4).Value = 44215
$ws.Cells.Item(63, 12).Value = 'Primera'
$ws.Cells.Item(63, 13).Value = 400
$ws.Cells.Item(63, 14).Value = 6000
$ws.Cells.Item(63, 15).Value = 6500
$ws.Cells.Item(63, 16).Value = 6250
$ws.Cells.Item(63, 18).Value = 'Perú'
$ws.Cells.Item(63, 19).Value = 1562

$ws.Cells.Item(64, 4).Value = 44432
$ws.Cells.Item(64, 12).Value = 'Primera'
$ws.Cells.Item(64, 13).Value = 200
$ws.Cells.Item(64, 14).Value = 8000
$ws.Cells.Item(64, 15).Value = 8500
$ws.Cells.Item(64, 16).Value = 8250
$ws.Cells.Item(64, 18).Value = 'Brasil'
$ws.Cells.Item(64, 19).Value = 2062

$ws.Cells.Item(65, 4).Value = 44600
$ws.Cells.Item(65, 12).Value = 'Primera'
$ws.Cells.Item(65, 13).Value = 200
$ws.Cells.Item(65, 14).Value = 7000
$ws.Cells.Item(65, 15).Value = 8000
$ws.Cells.Item(65, 16).Value = 7500
$ws.Cells.Item(65, 18).Value = 'Perú'
$ws.Cells.Item(65, 19).Value = 1875

$ws.Cells.Item(66, 4).Value = 44343
$ws.Cells.Item(66, 12).Value = 'Primera'
$ws.Cells.Item(66, 13).Value = 200
$ws.Cells.Item(66, 14).Value = 8000
$ws.Cells.Item(66, 15).Value = 9000
$ws.Cells.Item(66, 16).Value = 8500
$ws.Cells.Item(66, 18).Value = 'Perú'
$ws.Cells.Item(66, 19).Value = 2125

$ws.Cells.Item(67, 4).Value = 44721
$ws.Cells.Item(67, 12).Value = 'Primera'
$ws.Cells.Item(67, 13).Value = 180
$ws.Cells.Item(67, 14).Value = 9000
$ws.Cells.Item(67, 15).Value = 10000
$ws.Cells.Item(67, 16).Value = 9444
$ws.Cells.Item(67, 18).Value = 'Ecuador'
$ws.Cells.Item(67, 19).Value = 2361

$ws.Cells.Item(68, 4).Value = 44364
$ws.Cells.Item(68, 12).Value = 'Primera'
$ws.Cells.Item(68, 13).Value = 200
$ws.Cells.Item(68, 14).Value = 6500
$ws.Cells.Item(68, 15).Value = 7000
$ws.Cells.Item(68, 16).Value = 6750
$ws.Cells.Item(68, 18).Value = 'Perú'
$ws.Cells.Item(68, 19).Value = 1688

$ws.Cells.Item(69, 4).Value = 44222
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 200
$ws.Cells.Item(69, 14).Value = 6500
$ws.Cells.Item(69, 15).Value = 7000
$ws.Cells.Item(69, 16).Value = 6750
$ws.Cells.Item(69, 18).Value = 'Perú'
$ws.Cells.Item(69, 19).Value = 1688

$ws.Cells.Item(70, 4).Value = 44544
$ws.Cells.Item(70, 12).Value = 'Primera'
$ws.Cells.Item(70, 13).Value = 200
$ws.Cells.Item(70, 14).Value = 6000
$ws.Cells.Item(70, 15).Value = 6500
$ws.Cells.Item(70, 16).Value = 6250
$ws.Cells.Item(70, 18).Value = 'Perú'
$ws.Cells.Item(70, 19).Value = 1562

$ws.Cells.Item(71, 4).Value = 44749
$ws.Cells.Item(71, 12).Value = 'Primera'
$ws.Cells.Item(71, 13).Value = 200
$ws.Cells.Item(71, 14).Value = 7500
$ws.Cells.Item(71, 15).Value = 8000
$ws.Cells.Item(71, 16).Value = 7750
$ws.Cells.Item(71, 18).Value = 'Brasil'
$ws.Cells.Item(71, 19).Value = 1938

$ws.Cells.Item(72, 4).Value = 44161
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 200
$ws.Cells.Item(72, 14).Value = 7000
$ws.Cells.Item(72, 15).Value = 7500
$ws.Cells.Item(72, 16).Value = 7250
$ws.Cells.Item(72, 18).Value = 'Perú'
$ws.Cells.Item(72, 19).Value = 1812

$ws.Cells.Item(73, 4).Value = 44267
$ws.Cells.Item(73, 12).Value = 'Primera'
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 6000
$ws.Cells.Item(73, 15).Value = 6500
$ws.Cells.Item(73, 16).Value = 6250
$ws.Cells.Item(73, 18).Value = 'Perú'
$ws.Cells.Item(73, 19).Value = 1562

$ws.Cells.Item(74, 4).Value = 44678
$ws.Cells.Item(74, 12).Value = 'Primera'
$ws.Cells.Item(74, 13).Value = 180
$ws.Cells.Item(74, 14).Value = 7000
$ws.Cells.Item(74, 15).Value = 7500
$ws.Cells.Item(74, 16).Value = 7278
$ws.Cells.Item(74, 18).Value = 'Ecuador'
$ws.Cells.Item(74, 19).Value = 1820

$ws.Cells.Item(75, 4).Value = 44614
$ws.Cells.Item(75, 12).Value = 'Primera'
$ws.Cells.Item(75, 13).Value = 200
$ws.Cells.Item(75, 14).Value = 7000
$ws.Cells.Item(75, 15).Value = 7500
$ws.Cells.Item(75, 16).Value = 7250
$ws.Cells.Item(75, 18).Value = 'Perú'
$ws.Cells.Item(75, 19).Value = 1812

$ws.Cells.Item(76, 4).Value = 44239
$ws.Cells.Item(76, 12).Value = 'Primera'
$ws.Cells.Item(76, 13).Value = 200
$ws.Cells.Item(76, 14).Value = 6500
$ws.Cells.Item(76, 15).Value = 7000
$ws.Cells.Item(76, 16).Value = 6750
$ws.Cells.Item(76, 18).Value = 'Perú'
$ws.Cells.Item(76, 19).Value = 1688

$ws.Cells.Item(77, 4).Value = 44638
$ws.Cells.Item(77, 12).Value = 'Primera'
$ws.Cells.Item(77, 13).Value = 220
$ws.Cells.Item(77, 14).Value = 7000
$ws.Cells.Item(77, 15).Value = 7500
$ws.Cells.Item(77, 16).Value = 7227
$ws.Cells.Item(77, 18).Value = 'Ecuador'
$ws.Cells.Item(77, 19).Value = 1807

$ws.Cells.Item(78, 4).Value = 44194
$ws.Cells.Item(78, 12).Value = 'Primera'
$ws.Cells.Item(78, 13).Value = 400
$ws.Cells.Item(78, 14).Value = 6000
$ws.Cells.Item(78, 15).Value = 6500
$ws.Cells.Item(78, 16).Value = 6250
$ws.Cells.Item(78, 18).Value = 'Perú'
$ws.Cells.Item(78, 19).Value = 1562

$ws.Cells.Item(79, 4).Value = 44313
$ws.Cells.Item(79, 12).Value = 'Primera'
$ws.Cells.Item(79, 13).Value = 200
$ws.Cells.Item(79, 14).Value = 9000
$ws.Cells.Item(79, 15).Value = 10000
$ws.Cells.Item(79, 16).Value = 9500
$ws.Cells.Item(79, 18).Value = 'Ecuador'
$ws.Cells.Item(79, 19).Value = 2375

$ws.Cells.Item(80, 4).Value = 44616
$ws.Cells.Item(80, 12).Value = 'Primera'
$ws.Cells.Item(80, 13).Value = 200
$ws.Cells.Item(80, 14).Value = 7000
$ws.Cells.Item(80, 15).Value = 8000
$ws.Cells.Item(80, 16).Value = 7500
$ws.Cells.Item(80, 18).Value = 'Perú'
$ws.Cells.Item(80, 19).Value = 1875

$ws.Cells.Item(81, 4).Value = 44327
$ws.Cells.Item(81, 12).Value = 'Primera'
$ws.Cells.Item(81, 13).Value = 200
$ws.Cells.Item(81, 14).Value = 10000
$ws.Cells.Item(81, 15).Value = 11000
$ws.Cells.Item(81, 16).Value = 10500
$ws.Cells.Item(81, 18).Value = 'Perú'
$ws.Cells.Item(81, 19).Value = 2625

$ws.Cells.Item(82, 4).Value = 44673
$ws.Cells.Item(82, 12).Value = 'Primera'
$ws.Cells.Item(82, 13).Value = 200
$ws.Cells.Item(82, 14).Value = 7000
$ws.Cells.Item(82, 15).Value = 7500
$ws.Cells.Item(82, 16).Value = 7250
$ws.Cells.Item(82, 18).Value = 'Ecuador'
$ws.Cells.Item(82, 19).Value = 1812

$ws.Cells.Item(83, 4).Value = 44460
$ws.Cells.Item(83, 12).Value = 'Primera'
$ws.Cells.Item(83, 13).Value = 200
$ws.Cells.Item(83, 14).Value = 8000
$ws.Cells.Item(83, 15).Value = 8500
$ws.Cells.Item(83, 16).Value = 8250
$ws.Cells.Item(83, 18).Value = 'Brasil'
$ws.Cells.Item(83, 19).Value = 2062

$ws.Cells.Item(84, 4).Value = 44607
$ws.Cells.Item(84, 12).Value = 'Primera'
$ws.Cells.Item(84, 13).Value = 180
$ws.Cells.Item(84, 14).Value = 7000
$ws.Cells.Item(84, 15).Value = 7500
$ws.Cells.Item(84, 16).Value = 7222
$ws.Cells.Item(84, 18).Value = 'Perú'
$ws.Cells.Item(84, 19).Value = 1806

$ws.Cells.Item(85, 4).Value = 44292
$ws.Cells.Item(85, 12).Value = 'Primera'
$ws.Cells.Item(85, 13).Value = 400
$ws.Cells.Item(85, 14).Value = 7000
$ws.Cells.Item(85, 15).Value = 7500
$ws.Cells.Item(85, 16).Value = 7250
$ws.Cells.Item(85, 18).Value = 'Perú'
$ws.Cells.Item(85, 19).Value = 1812

$ws.Cells.Item(86, 4).Value = 44348
$ws.Cells.Item(86, 12).Value = 'Primera'
$ws.Cells.Item(86, 13).Value = 200
$ws.Cells.Item(86, 14).Value = 8000
$ws.Cells.Item(86, 15).Value = 8500
$ws.Cells.Item(86, 16).Value = 8250
$ws.Cells.Item(86, 18).Value = 'Perú'
$ws.Cells.Item(86, 19).Value = 2062

$ws.Cells.Item(87, 4).Value = 44370
$ws.Cells.Item(87, 12).Value = 'Primera'
$ws.Cells.Item(87, 13).Value = 200
$ws.Cells.Item(87, 14).Value = 6000
$ws.Cells.Item(87, 15).Value = 6500
$ws.Cells.Item(87, 16).Value = 6250
$ws.Cells.Item(87, 18).Value = 'Perú'
$ws.Cells.Item(87, 19).Value = 1562

$ws.Cells.Item(88, 4).Value = 44421
$ws.Cells.Item(88, 12).Value = 'Primera'
$ws.Cells.Item(88, 13).Value = 200
$ws.Cells.Item(88, 14).Value = 8500
$ws.Cells.Item(88, 15).Value = 9000
$ws.Cells.Item(88, 16).Value = 8750
$ws.Cells.Item(88, 18).Value = 'Brasil'
$ws.Cells.Item(88, 19).Value = 2188

$ws.Cells.Item(89, 4).Value = 44622
$ws.Cells.Item(89, 12).Value = 'Primera'
$ws.Cells.Item(89, 13).Value = 150
$ws.Cells.Item(89, 14).Value = 7000
$ws.Cells.Item(89, 15).Value = 7500
$ws.Cells.Item(89, 16).Value = 7233
$ws.Cells.Item(89, 18).Value = 'Perú'
$ws.Cells.Item(89, 19).Value = 1808

$ws.Cells.Item(90, 4).Value = 44565
$ws.Cells.Item(90, 12).Value = 'Primera'
$ws.Cells.Item(90, 13).Value = 140
$ws.Cells.Item(90, 14).Value = 6500
$ws.Cells.Item(90, 15).Value = 7000
$ws.Cells.Item(90, 16).Value = 6786
$ws.Cells.Item(90, 18).Value = 'Perú'
$ws.Cells.Item(90, 19).Value = 1696

$ws.Cells.Item(91, 4).Value = 44259
$ws.Cells.Item(91, 12).Value = 'Primera'
$ws.Cells.Item(91, 13).Value = 400
$ws.Cells.Item(91, 14).Value = 6000
$ws.Cells.Item(91, 15).Value = 6500
$ws.Cells.Item(91, 16).Value = 6250
$ws.Cells.Item(91, 18).Value = 'Perú'
$ws.Cells.Item(91, 19).Value = 1562

$ws.Cells.Item(92, 4).Value = 44202
$ws.Cells.Item(92, 12).Value = 'Primera'
$ws.Cells.Item(92, 13).Value = 200
$ws.Cells.Item(92, 14).Value = 6000
$ws.Cells.Item(92, 15).Value = 6500
$ws.Cells.Item(92, 16).Value = 6250
$ws.Cells.Item(92, 18).Value = 'Perú'
$ws.Cells.Item(92, 19).Value = 1562

$ws.Cells.Item(93, 4).Value = 44371
$ws.Cells.Item(93, 12).Value = 'Primera'
$ws.Cells.Item(93, 13).Value = 200
$ws.Cells.Item(93, 14).Value = 6500
$ws.Cells.Item(93, 15).Value = 7000
$ws.Cells.Item(93, 16).Value = 6750
$ws.Cells.Item(93, 18).Value = 'Perú'
$ws.Cells.Item(93, 19).Value = 1688

$ws.Cells.Item(94, 4).Value = 44162
$ws.Cells.Item(94, 12).Value = 'Primera'
$ws.Cells.Item(94, 13).Value = 200
$ws.Cells.Item(94, 14).Value = 7000
$ws.Cells.Item(94, 15).Value = 7500
$ws.Cells.Item(94, 16).Value = 7250
$ws.Cells.Item(94, 18).Value = 'Perú'
$ws.Cells.Item(94, 19).Value = 1812

$ws.Cells.Item(95, 4).Value = 44726
$ws.Cells.Item(95, 12).Value = 'Primera'
$ws.Cells.Item(95, 13).Value = 200
$ws.Cells.Item(95, 14).Value = 8000
$ws.Cells.Item(95, 15).Value = 9000
$ws.Cells.Item(95, 16).Value = 8500
$ws.Cells.Item(95, 18).Value = 'Brasil'
$ws.Cells.Item(95, 19).Value = 2125

$ws.Cells.Item(96, 4).Value = 44427
$ws.Cells.Item(96, 12).Value = 'Primera'
$ws.Cells.Item(96, 13).Value = 200
$ws.Cells.Item(96, 14).Value = 8500
$ws.Cells.Item(96, 15).Value = 9000
$ws.Cells.Item(96, 16).Value = 8750
$ws.Cells.Item(96, 18).Value = 'Brasil'
$ws.Cells.Item(96, 19).Value = 2188

$ws.Cells.Item(97, 4).Value = 44441
$ws.Cells.Item(97, 12).Value = 'Primera'
$ws.Cells.Item(97, 13).Value = 200
$ws.Cells.Item(97, 14).Value = 7500
$ws.Cells.Item(97, 15).Value = 8000
$ws.Cells.Item(97, 16).Value = 7750
$ws.Cells.Item(97, 18).Value = 'Brasil'
$ws.Cells.Item(97, 19).Value = 1938

$ws.Cells.Item(98, 4).Value = 44237
$ws.Cells.Item(98, 12).Value = 'Primera'
$ws.Cells.Item(98, 13).Value = 200
$ws.Cells.Item(98, 14).Value = 6000
$ws.Cells.Item(98, 15).Value = 6500
$ws.Cells.Item(98, 16).Value = 6250
$ws.Cells.Item(98, 18).Value = 'Perú'
$ws.Cells.Item(98, 19).Value = 1562

$ws.Cells.Item(99, 4).Value = 44628
$ws.Cells.Item(99, 12).Value = 'Primera'
$ws.Cells.Item(99, 13).Value = 200
$ws.Cells.Item(99, 14).Value = 7500
$ws.Cells.Item(99, 15).Value = 8000
$ws.Cells.Item(99, 16).Value = 7750
$ws.Cells.Item(99, 18).Value = 'Ecuador'
$ws.Cells.Item(99, 19).Value = 1938

$ws.Cells.Item(100, 4).Value = 44483
$ws.Cells.Item(100, 12).Value = 'Primera'
$ws.Cells.Item(100, 13).Value = 200
$ws.Cells.Item(100, 14).Value = 7000
$ws.Cells.Item(100, 15).Value = 7500
$ws.Cells.Item(100, 16).Value = 7250
$ws.Cells.Item(100, 18).Value = 'Perú'
$ws.Cells.Item(100, 19).Value = 1812

$ws.Cells.Item(101, 4).Value = 44253
$ws.Cells.Item(101, 12).Value = 'Primera'
$ws.Cells.Item(101, 13).Value = 400
$ws.Cells.Item(101, 14).Value = 6000
$ws.Cells.Item(101, 15).Value = 6500
$ws.Cells.Item(101, 16).Value = 6250
$ws.Cells.Item(101, 18).Value = 'Perú'
$ws.Cells.Item(101, 19).Value = 1562

$ws.Cells.Item(102, 4).Value = 44204
$ws.Cells.Item(102, 12).Value = 'Primera'
$ws.Cells.Item(102, 13).Value = 200
$ws.Cells.Item(102, 14).Value = 5000
$ws.Cells.Item(102, 15).Value = 5500
$ws.Cells.Item(102, 16).Value = 5250
$ws.Cells.Item(102, 18).Value = 'Perú'
$ws.Cells.Item(102, 19).Value = 1312

$ws.Cells.Item(103, 4).Value = 44266
$ws.Cells.Item(103, 12).Value = 'Primera'
$ws.Cells.Item(103, 13).Value = 200
$ws.Cells.Item(103, 14).Value = 6000
$ws.Cells.Item(103, 15).Value = 6500
$ws.Cells.Item(103, 16).Value = 6250
$ws.Cells.Item(103, 18).Value = 'Perú'
$ws.Cells.Item(103, 19).Value = 1562

$ws.Cells.Item(104, 4).Value = 44229
$ws.Cells.Item(104, 12).Value = 'Primera'
$ws.Cells.Item(104, 13).Value = 400
$ws.Cells.Item(104, 14).Value = 6500
$ws.Cells.Item(104, 15).Value = 7000
$ws.Cells.Item(104, 16).Value = 6750
$ws.Cells.Item(104, 18).Value = 'Perú'
$ws.Cells.Item(104, 19).Value = 1688

$ws.Cells.Item(105, 4).Value = 44231
$ws.Cells.Item(105, 12).Value = 'Primera'
$ws.Cells.Item(105, 13).Value = 400
$ws.Cells.Item(105, 14).Value = 6000
$ws.Cells.Item(105, 15).Value = 6500
$ws.Cells.Item(105, 16).Value = 6250
$ws.Cells.Item(105, 18).Value = 'Perú'
$ws.Cells.Item(105, 19).Value = 1562

$ws.Cells.Item(106, 4).Value = 44336
$ws.Cells.Item(106, 12).Value = 'Primera'
$ws.Cells.Item(106, 13).Value = 200
$ws.Cells.Item(106, 14).Value = 10000
$ws.Cells.Item(106, 15).Value = 11000
$ws.Cells.Item(106, 16).Value = 10500
$ws.Cells.Item(106, 18).Value = 'Perú'
$ws.Cells.Item(106, 19).Value = 2625

$ws.Cells.Item(107, 4).Value = 44488
$ws.Cells.Item(107, 12).Value = 'Primera'
$ws.Cells.Item(107, 13).Value = 280
$ws.Cells.Item(107, 14).Value = 8000
$ws.Cells.Item(107, 15).Value = 9000
$ws.Cells.Item(107, 16).Value = 8464
$ws.Cells.Item(107, 18).Value = 'Perú'
$ws.Cells.Item(107, 19).Value = 2116

$ws.Cells.Item(108, 4).Value = 44196
$ws.Cells.Item(108, 12).Value = 'Primera'
$ws.Cells.Item(108, 13).Value = 200
$ws.Cells.Item(108, 14).Value = 6000
$ws.Cells.Item(108, 15).Value = 6500
$ws.Cells.Item(108, 16).Value = 6250
$ws.Cells.Item(108, 18).Value = 'Perú'
$ws.Cells.Item(108, 19).Value = 1562

$ws.Cells.Item(109, 4).Value = 44596
$ws.Cells.Item(109, 12).Value = 'Primera'
$ws.Cells.Item(109, 13).Value = 170
$ws.Cells.Item(109, 14).Value = 6500
$ws.Cells.Item(109, 15).Value = 7000
$ws.Cells.Item(109, 16).Value = 6735
$ws.Cells.Item(109, 18).Value = 'Ecuador'
$ws.Cells.Item(109, 19).Value = 1684

$ws.Cells.Item(110, 4).Value = 44496
$ws.Cells.Item(110, 12).Value = 'Primera'
$ws.Cells.Item(110, 13).Value = 200
$ws.Cells.Item(110, 14).Value = 7500
$ws.Cells.Item(110, 15).Value = 8000
$ws.Cells.Item(110, 16).Value = 7750
$ws.Cells.Item(110, 18).Value = 'Perú'
$ws.Cells.Item(110, 19).Value = 1938

$ws.Cells.Item(111, 4).Value = 44496
$ws.Cells.Item(111, 12).Value = 'Segunda'
$ws.Cells.Item(111, 13).Value = 100
$ws.Cells.Item(111, 14).Value = 7000
$ws.Cells.Item(111, 15).Value = 7000
$ws.Cells.Item(111, 16).Value = 7000
$ws.Cells.Item(111, 18).Value = 'Perú'
$ws.Cells.Item(111, 19).Value = 1750

$ws.Cells.Item(112, 4).Value = 44399
$ws.Cells.Item(112, 12).Value = 'Primera'
$ws.Cells.Item(112, 13).Value = 200
$ws.Cells.Item(112, 14).Value = 7000
$ws.Cells.Item(112, 15).Value = 7500
$ws.Cells.Item(112, 16).Value = 7250
$ws.Cells.Item(112, 18).Value = 'Brasil'
$ws.Cells.Item(112, 19).Value = 1812

$ws.Cells.Item(113, 4).Value = 44497
$ws.Cells.Item(113, 12).Value = 'Primera'
$ws.Cells.Item(113, 13).Value = 200
$ws.Cells.Item(113, 14).Value = 7000
$ws.Cells.Item(113, 15).Value = 7500
$ws.Cells.Item(113, 16).Value = 7250
$ws.Cells.Item(113, 18).Value = 'Perú'
$ws.Cells.Item(113, 19).Value = 1812

$ws.Cells.Item(114, 4).Value = 44285
$ws.Cells.Item(114, 12).Value = 'Primera'
$ws.Cells.Item(114, 13).Value = 200
$ws.Cells.Item(114, 14).Value = 6000
$ws.Cells.Item(114, 15).Value = 6500
$ws.Cells.Item(114, 16).Value = 6250
$ws.Cells.Item(114, 18).Value = 'Perú'
$ws.Cells.Item(114, 19).Value = 1562

$ws.Cells.Item(115, 4).Value = 44657
$ws.Cells.Item(115, 12).Value = 'Primera'
$ws.Cells.Item(115, 13).Value = 200
$ws.Cells.Item(115, 14).Value = 7500
$ws.Cells.Item(115, 15).Value = 8000
$ws.Cells.Item(115, 16).Value = 7750
$ws.Cells.Item(115, 18).Value = 'Perú'
$ws.Cells.Item(115, 19).Value = 1938

$ws.Cells.Item(116, 4).Value = 44391
$ws.Cells.Item(116, 12).Value = 'Primera'
$ws.Cells.Item(116, 13).Value = 200
$ws.Cells.Item(116, 14).Value = 6000
$ws.Cells.Item(116, 15).Value = 6500
$ws.Cells.Item(116, 16).Value = 6250
$ws.Cells.Item(116, 18).Value = 'Brasil'
$ws.Cells.Item(116, 19).Value = 1562

$ws.Cells.Item(117, 4).Value = 44453
$ws.Cells.Item(117, 12).Value = 'Primera'
$ws.Cells.Item(117, 13).Value = 200
$ws.Cells.Item(117, 14).Value = 8000
$ws.Cells.Item(117, 15).Value = 8500
$ws.Cells.Item(117, 16).Value = 8250
$ws.Cells.Item(117, 18).Value = 'Brasil'
$ws.Cells.Item(117, 19).Value = 2062

$ws.Cells.Item(118, 4).Value = 44572
$ws.Cells.Item(118, 12).Value = 'Primera'
$ws.Cells.Item(118, 13).Value = 200
$ws.Cells.Item(118, 14).Value = 6000
$ws.Cells.Item(118, 15).Value = 6500
$ws.Cells.Item(118, 16).Value = 6250
$ws.Cells.Item(118, 18).Value = 'Perú'
$ws.Cells.Item(118, 19).Value = 1562

Write-Output "done"